$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B124').Value = 6937238
$ws.Range('E124').Value = 'PAOK Salonika'
$ws.Range('F124').Value = 'Giannina'
$ws.Range('G124').Value = 4
$ws.Range('H124').Value = 0
$ws.Range('J124').Value = 0
$ws.Range('K124').Value = 'H'
$ws.Range('L124').Value = 1.111
$ws.Range('M124').Value = 9
$ws.Range('N124').Value = 23
$ws.Range('O124').Value = 1.25
$ws.Range('P124').Value = 6
$ws.Range('Q124').Value = 9
$ws.Range('R124').Value = -1.75
$ws.Range('S124').Value = 2.025
$ws.Range('T124').Value = 1.825
$ws.Range('U124').Value = 2.75
$ws.Range('X124').Value = 0.25
$ws.Range('Y124').Value = -1
$ws.Range('AA124').Value = 1.025
$ws.Range('AB124').Value = -1
$ws.Range('B125').Value = 6936857
$ws.Range('E125').Value = 'AEK Athens'
$ws.Range('F125').Value = 'Panathinaikos'
$ws.Range('G125').Value = 2
$ws.Range('H125').Value = 2
$ws.Range('J125').Value = 1
$ws.Range('K125').Value = 'D'
$ws.Range('L125').Value = 1.909
$ws.Range('M125').Value = 3.5
$ws.Range('N125').Value = 4.2
$ws.Range('O125').Value = 2.15
$ws.Range('P125').Value = 3.2
$ws.Range('Q125').Value = 3.5
$ws.Range('R125').Value = -0.25
$ws.Range('S125').Value = 1.85
$ws.Range('T125').Value = 2
$ws.Range('U125').Value = 2
$ws.Range('X125').Value = -1
$ws.Range('Y125').Value = 2.2
$ws.Range('AA125').Value = -0.5
$ws.Range('AB125').Value = 0.5
$ws.Range('B200').Value = 7920453
$ws.Range('E200').Value = 'Panetolikos'
$ws.Range('F200').Value = 'Volos NFC'
$ws.Range('G200').Value = 0
$ws.Range('H200').Value = 1
$ws.Range('I200').Value = 0
$ws.Range('J200').Value = 1
$ws.Range('L200').Value = 2.3
$ws.Range('M200').Value = 3
$ws.Range('N200').Value = 3.4
$ws.Range('P200').Value = 3.1
$ws.Range('Q200').Value = 3.8
$ws.Range('S200').Value = 1.8
$ws.Range('T200').Value = 2.05
$ws.Range('U200').Value = 2.25
$ws.Range('Z200').Value = 2.8
$ws.Range('AB200').Value = 1.05
$ws.Range('AC200').Value = -1
$ws.Range('AD200').Value = 0.825
$ws.Range('B201').Value = 7920450
$ws.Range('E201').Value = 'Asteras Tripolis'
$ws.Range('F201').Value = 'Kifisias FC'
$ws.Range('G201').Value = 1
$ws.Range('H201').Value = 2
$ws.Range('I201').Value = 1
$ws.Range('J201').Value = 0
$ws.Range('L201').Value = 2.05
$ws.Range('M201').Value = 3.3
$ws.Range('N201').Value = 3.6
$ws.Range('P201').Value = 3.5
$ws.Range('Q201').Value = 3.4
$ws.Range('S201').Value = 1.85
$ws.Range('T201').Value = 2
$ws.Range('U201').Value = 2.75
$ws.Range('Z201').Value = 2.4
$ws.Range('AB201').Value = 1
$ws.Range('AC201').Value = 0.5125
$ws.Range('AD201').Value = -0.5
